# Feria Lagunitas de Puerto Montt - Pepino dulce
# Weekly update: insert a new price record (row 79) for the latest week,
# pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 79, shifting rows 79:89 -> 80:90
$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with this week's data
$ws.Cells.Item(79, 1).Value = 4
$ws.Cells.Item(79, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value = "Los Lagos"
$ws.Cells.Item(79, 4).Value = 45034
$ws.Cells.Item(79, 5).Value = 10
$ws.Cells.Item(79, 6).Value = 100112043
$ws.Cells.Item(79, 7).Value = "Pepino dulce"
$ws.Cells.Item(79, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 80
$ws.Cells.Item(79, 11).Value = 19000
$ws.Cells.Item(79, 12).Value = 19000
$ws.Cells.Item(79, 13).Value = 19000
$ws.Cells.Item(79, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(79, 16).Value = 1056
$ws.Cells.Item(79, 17).Value = 18
$ws.Cells.Item(79, 18).Value = "Hortaliza"
